# Apply the changes described by the diff:
#  - PersonalDetails: rename/retext CoreCompetency bullet values
#  - ProjectsDetails: add a ProjectYear column, and swap the
#    "Forecasting Drug Sales" / "Editors are Gatekeepers of Science"
#    project rows (4 and 5)
#  - Update sheet selections / active sheet to match final saved state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. PersonalDetails - CoreCompetency wording updates
# ---------------------------------------------------------------
$wsPersonal = $wb.Worksheets.Item("PersonalDetails")
$wsPersonal.Range("B10").Value = "Domain knowledge in Finance and Innovation"
$wsPersonal.Range("B11").Value = "Data driven decision making"
$wsPersonal.Range("B9").Value = "Regression and classification techniques"

# ---------------------------------------------------------------
# 2. ProjectsDetails - add ProjectYear column + swap rows 5 & 6
# ---------------------------------------------------------------
$wsProjects = $wb.Worksheets.Item("ProjectsDetails")

# New ProjectYear column
$wsProjects.Range("E1").Value = "ProjectYear"
$wsProjects.Range("E2").Value = 2019
$wsProjects.Range("E3").Value = 2019
$wsProjects.Range("E4").Value = 2020

# Swap content of rows 5 (ProjectId 4) and 6 (ProjectId 5)
$wsProjects.Range("B5").Value = "Editors are Gatekeepers of Science"
$wsProjects.Range("C5").Value = "Editors"
$wsProjects.Range("B6").Value = "Forecasting Drug Sales"
$wsProjects.Range("C6").Value = "DrugSales"

# ---------------------------------------------------------------
# 3. Sheet view / selection / active-sheet bookkeeping
# ---------------------------------------------------------------
$wsPersonal.Activate()
$wsPersonal.Range("B9").Select() | Out-Null

$wsEducation = $wb.Worksheets.Item("EducationDetails")
$wsWorkExp = $wb.Worksheets.Item("WorkExperienceDetails")
$wsWorkBullets = $wb.Worksheets.Item("WorkExperienceBullets")
$wsProjectBullets = $wb.Worksheets.Item("ProjectsBullets")

$wsPersonal.Range("B17").Select() | Out-Null
$wsProjects.Range("E5").Select() | Out-Null
$wsProjectBullets.Range("J31").Select() | Out-Null

$wsProjects.Activate()
